$wb = $excel.ActiveWorkbook

# --- Incidental selection-state changes on the other sheets, made before
#     switching over to the new sheet ---
$germany = $wb.Worksheets.Item("Germany")
[void]$germany.Range("A12").Select()

$slovakia = $wb.Worksheets.Item("Slovakia")
[void]$slovakia.Cells.Select()

# --- Create the new "Italy" sheet as a copy of "Slovakia", placed at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Copy($null, $lastSheet)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"
$italy.Activate()

# Remove the two market-specific rows that don't apply (old rows 9 & 10:
# "XLM800-STI" / "XLM800-Zetfas"), shrinking the sheet from 12 to 10 rows.
$italy.Rows("9:10").Delete()

# Fill in the Italy-specific market name and reference code.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2454/T2453"

# B4 lost its bordered style when its content was retyped.
$italy.Range("B4").Borders.LineStyle = 0

# Leave the cursor where the author left it.
[void]$italy.Range("C16").Select()
